# 2024大作业自评表.xlsx — "update part7 of readme"
# The author re-reviewed the self-assessment checklist and marked several
# more requirement rows as completed ("是"), then left the selection near
# the bottom of the newly-updated block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Part 7-ish requirements that got checked off as done.
$ws.Range("B26").Value = "是"
$ws.Range("B34").Value = "是"
$ws.Range("B38").Value = "是"
$ws.Range("B39").Value = "是"
$ws.Range("B40").Value = "是"

# Leave the cursor where the author was last working.
$ws.Range("C41").Select() | Out-Null

# The window background (lt1) theme colour reverted from the custom
# light-green tint back to plain white.
$themeColors = $wb.Theme.ThemeColorScheme
$background1 = $themeColors.Colors(2)
$background1.RGB = 16777215
